# Update crypto price/volume data (cryptos list refresh) and fix an
# Avalanche / WrappedliquidstakedEther2.0 row-order swap (rows 14 & 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.434.16"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.523.15"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "3.522.48"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.120.22"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "3.516.52"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "67.423.78"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +7.64%  "
$ws.Range("D26").Value = "3.665.54"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.49%  "
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("E33").Value = "  +4.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").Value = "3.515.54"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.69%  "
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -4.46%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E51").Value = "  -3.08%  "
